# Auto-generated Excel COM-interop script
# Applies numeric value updates (simulated scheduled market-price refresh)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 37500
$ws.Range("J63").Value = 37500
$ws.Range("L63").Value = 37500
$ws.Range("N63").Value = -38748
$ws.Range("H66").Value = 37500
$ws.Range("J66").Value = 37500
$ws.Range("L66").Value = 112500
$ws.Range("N66").Value = -118740
$ws.Range("H100").Value = 2332.1428
$ws.Range("I100").Value = 1934.4445
$ws.Range("J100").Value = 2630.4167
$ws.Range("K100").Value = 1934.4445
$ws.Range("L100").Value = 2630.4167
$ws.Range("M100").Value = -1393.4445
$ws.Range("N100").Value = -3712.4167
$ws.Range("H107").Value = 1083
$ws.Range("I107").Value = 966.6667
$ws.Range("K107").Value = 966.6667
$ws.Range("M107").Value = 953.3333
$ws.Range("H112").Value = 2623.158
$ws.Range("J112").Value = 2713.3333
$ws.Range("L112").Value = 8139.999899999999
$ws.Range("N112").Value = -10355.9999
$ws.Range("H129").Value = 1021.2353
$ws.Range("I129").Value = 1463.091
$ws.Range("J129").Value = 955.5540999999999
$ws.Range("K129").Value = 4389.272999999999
$ws.Range("L129").Value = 2866.6623
$ws.Range("M129").Value = 610.7270000000008
$ws.Range("N129").Value = -12866.6623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3816.111
$ws.Range("J63").Value = 5345
$ws.Range("L63").Value = 5345
$ws.Range("N63").Value = -6717
$ws.Range("H66").Value = 3816.111
$ws.Range("J66").Value = 5345
$ws.Range("L66").Value = 26725
$ws.Range("N66").Value = -33589
$ws.Range("H110").Value = 1385.3478
$ws.Range("I110").Value = 1380.1364
$ws.Range("K110").Value = 1380.1364
$ws.Range("M110").Value = 664.8635999999999
$ws.Range("H122").Value = 2261.1667
$ws.Range("I122").Value = 2127.1667
$ws.Range("J122").Value = 2395.1667
$ws.Range("K122").Value = 6381.500100000001
$ws.Range("L122").Value = 7185.500100000001
$ws.Range("M122").Value = -3931.500100000001
$ws.Range("N122").Value = -12085.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 35486.668
$ws.Range("J35").Value = 35486.668
$ws.Range("L35").Value = 35486.668
$ws.Range("N35").Value = -36106.668
$ws.Range("H82").Value = 15193.462
$ws.Range("H85").Value = 15193.462
$ws.Range("H105").Value = 2920.4546
$ws.Range("I105").Value = 1848.25
$ws.Range("K105").Value = 1848.25
$ws.Range("M105").Value = -101.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 4692.3076
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 4692.3076
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 4692.3076
$ws.Range("M6").ClearContents()  # was -37
$ws.Range("N6").Value = -4918.3076
$ws.Range("H31").Value = 6541823.5
$ws.Range("I31").Value = 3507.1428
$ws.Range("J31").Value = 9015781
$ws.Range("K31").Value = 3507.1428
$ws.Range("L31").Value = 9015781
$ws.Range("M31").Value = -3212.1428
$ws.Range("N31").Value = -9016371
$ws.Range("H34").Value = 6541823.5
$ws.Range("I34").Value = 3507.1428
$ws.Range("J34").Value = 9015781
$ws.Range("K34").Value = 3507.1428
$ws.Range("L34").Value = 9015781
$ws.Range("M34").Value = -3305.1428
$ws.Range("N34").Value = -9016185
$ws.Range("H41").Value = 28964
$ws.Range("J41").Value = 28964
$ws.Range("L41").Value = 28964
$ws.Range("N41").Value = -29820
$ws.Range("H50").Value = 34490
$ws.Range("J50").Value = 34490
$ws.Range("L50").Value = 34490
$ws.Range("N50").Value = -35740
$ws.Range("H51").Value = 100031064
$ws.Range("J51").Value = 38830
$ws.Range("L51").Value = 38830
$ws.Range("N51").Value = -40302
$ws.Range("H59").Value = 29979
$ws.Range("J59").Value = 27954
$ws.Range("L59").Value = 27954
$ws.Range("N59").Value = -30244
$ws.Range("H60").Value = 11331.286
$ws.Range("I60").Value = 7000
$ws.Range("J60").Value = 11547.85
$ws.Range("K60").Value = 7000
$ws.Range("L60").Value = 11547.85
$ws.Range("M60").Value = -6489
$ws.Range("N60").Value = -12569.85
$ws.Range("H61").Value = 100031064
$ws.Range("J61").Value = 38830
$ws.Range("L61").Value = 38830
$ws.Range("N61").Value = -39526
$ws.Range("H74").Value = 13000
$ws.Range("J74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("N74").Value = -14748
$ws.Range("H77").Value = 13000
$ws.Range("J77").Value = 13000
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -47736

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 8607.291999999999
$ws.Range("I107").Value = 11376.889
$ws.Range("J107").Value = 6945.533
$ws.Range("K107").Value = 34130.667
$ws.Range("L107").Value = 20836.599
$ws.Range("M107").Value = -32210.667
$ws.Range("N107").Value = -24676.599
$ws.Range("H112").Value = 12503693
$ws.Range("I112").Value = 100000960
$ws.Range("J112").Value = 4082.8572
$ws.Range("K112").Value = 300002880
$ws.Range("L112").Value = 12248.5716
$ws.Range("M112").Value = -300001772
$ws.Range("N112").Value = -14464.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5242.067
$ws.Range("I97").Value = 3510.9092
$ws.Range("J97").Value = 10002.75
$ws.Range("K97").Value = 3510.9092
$ws.Range("L97").Value = 10002.75
$ws.Range("M97").Value = -3014.9092
$ws.Range("N97").Value = -10994.75
$ws.Range("H102").Value = 1300
$ws.Range("I102").Value = 1325
$ws.Range("J102").Value = 1200
$ws.Range("K102").Value = 1325
$ws.Range("L102").Value = 1200
$ws.Range("M102").Value = 297
$ws.Range("N102").Value = -4444
$ws.Range("H122").Value = 1953
$ws.Range("J122").Value = 2466.6667
$ws.Range("L122").Value = 7400.000100000001
$ws.Range("N122").Value = -12300.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3518.2856
$ws.Range("I40").Value = 3433.4443
$ws.Range("J40").Value = 3671
$ws.Range("K40").Value = 3433.4443
$ws.Range("L40").Value = 3671
$ws.Range("M40").Value = -3297.4443
$ws.Range("N40").Value = -3943
$ws.Range("H61").Value = 2315.9
$ws.Range("I61").Value = 2495.4443
$ws.Range("J61").Value = 700
$ws.Range("K61").Value = 2495.4443
$ws.Range("L61").Value = 700
$ws.Range("M61").Value = -2293.4443
$ws.Range("N61").Value = -1104
$ws.Range("H113").Value = 2315.9
$ws.Range("I113").Value = 2495.4443
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 2495.4443
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = -325.4443000000001
$ws.Range("N113").Value = -5040
$ws.Range("H122").Value = 93218.55
$ws.Range("I122").Value = 93218.55
$ws.Range("K122").Value = 279655.65
$ws.Range("M122").Value = -277205.65

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 550
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 1650
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 800
$ws.Range("N122").Value = -7900
